$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 276; existing rows 276-337 shift down to 277-338.
$ws.Rows.Item(276).Insert()

# Populate the new row 276 with the new data record.
$ws.Range("A276").Value = 5
$ws.Range("B276").Value = "Macroferia Regional de Talca"
$ws.Range("C276").Value = "Maule"
$ws.Range("D276").Value = 44782
$ws.Range("E276").Value = 7
$ws.Range("F276").Value = 100112003
$ws.Range("G276").Value = "Ajo"
$ws.Range("H276").Value = "Chino"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 300
$ws.Range("K276").Value = 28000
$ws.Range("L276").Value = 28000
$ws.Range("M276").Value = 28000
$ws.Range("N276").Value = "$/malla 10 kilos"
$ws.Range("O276").Value = "China"
$ws.Range("P276").Value = 2800
$ws.Range("Q276").Value = 10
$ws.Range("R276").Value = "Hortaliza"
